$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 557
$ws.Cells.Item(3, 6).Value = 10384
$ws.Cells.Item(4, 6).Value = 228
$ws.Cells.Item(5, 6).Value = 95
$ws.Cells.Item(6, 6).Value = 6940
$ws.Cells.Item(7, 6).Value = 654
$ws.Cells.Item(8, 6).Value = 133
$ws.Cells.Item(9, 6).Value = 12087
$ws.Cells.Item(11, 6).Value = 1308
$ws.Cells.Item(12, 6).Value = 1278
$ws.Cells.Item(13, 6).Value = 5396
$ws.Cells.Item(14, 6).Value = 901
$ws.Cells.Item(16, 6).Value = 359
$ws.Cells.Item(17, 6).Value = 193
$ws.Cells.Item(18, 6).Value = 1420
$ws.Cells.Item(19, 6).Value = 336
$ws.Cells.Item(20, 6).Value = 1998
$ws.Cells.Item(21, 6).Value = 1017
$ws.Cells.Item(22, 6).Value = 1524
$ws.Cells.Item(25, 6).Value = 499
$ws.Cells.Item(26, 6).Value = 721
$ws.Cells.Item(27, 6).Value = 2967
$ws.Cells.Item(28, 6).Value = 249
$ws.Cells.Item(29, 6).Value = 2014
$ws.Cells.Item(30, 6).Value = 108
$ws.Cells.Item(31, 6).Value = 1675
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 124
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 93
$ws.Cells.Item(36, 6).Value = 3672
$ws.Cells.Item(37, 6).Value = 4367
$ws.Cells.Item(38, 6).Value = 267
$ws.Cells.Item(39, 6).Value = 125
$ws.Cells.Item(42, 6).Value = 575
$ws.Cells.Item(43, 6).Value = 35
$ws.Cells.Item(45, 6).Value = 285
$ws.Cells.Item(46, 6).Value = 24
$ws.Cells.Item(47, 6).Value = 26
$ws.Cells.Item(48, 6).Value = 4283
$ws.Cells.Item(49, 6).Value = 175

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 15
$ws.Cells.Item(4, 6).Value = 15
$ws.Cells.Item(5, 6).Value = 36
$ws.Cells.Item(9, 6).Value = 68
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(12, 6).Value = 75
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 11
$ws.Cells.Item(21, 6).Value = 5
$ws.Cells.Item(23, 6).Value = 74
$ws.Cells.Item(27, 6).Value = 49
$ws.Cells.Item(28, 6).Value = 2
$ws.Cells.Item(29, 6).Value = 61
$ws.Cells.Item(30, 6).Value = 5
$ws.Cells.Item(31, 6).Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 6448

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 557
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 654
$ws.Cells.Item(9, 6).Value = 12087
$ws.Cells.Item(10, 6).Value = 12684
$ws.Cells.Item(11, 6).Value = 35
$ws.Cells.Item(12, 6).Value = 1308
$ws.Cells.Item(14, 6).Value = 901
$ws.Cells.Item(19, 6).Value = 1017
$ws.Cells.Item(20, 6).Value = 1524
$ws.Cells.Item(21, 6).Value = 879
$ws.Cells.Item(22, 6).Value = 5
$ws.Cells.Item(23, 6).Value = 16
$ws.Cells.Item(24, 6).Value = 499
$ws.Cells.Item(25, 6).Value = 721
$ws.Cells.Item(28, 6).Value = 2014
$ws.Cells.Item(29, 6).Value = 6
$ws.Cells.Item(30, 6).Value = 1675
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 5
$ws.Cells.Item(34, 6).Value = 12
$ws.Cells.Item(36, 6).Value = 37
$ws.Cells.Item(37, 6).Value = 4367
$ws.Cells.Item(39, 6).Value = 267
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 576
$ws.Cells.Item(43, 6).Value = 913
$ws.Cells.Item(44, 6).Value = 285
$ws.Cells.Item(46, 6).Value = 4283
$ws.Cells.Item(47, 6).Value = 175
